$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Backlog sheet: add two new backlog items (rows 15-16) and fill in the
# feature names for the already-present empty rows (17-19).
# ---------------------------------------------------------------------------
$backlog = $wb.Worksheets.Item("Backlog")

$backlog.Range("A15").Value = 13
$backlog.Range("B15").Value = "Score / reeks"
$backlog.Range("D15").Value = "M"
$backlog.Range("F15").Value = "x"

$backlog.Range("A16").Value = 14
$backlog.Range("B16").Value = "Levens / kansen"
$backlog.Range("D16").Value = "M"
$backlog.Range("F16").Value = "x"

$backlog.Range("B17").Value = "Ranking"
$backlog.Range("B18").Value = "Inlog"
$backlog.Range("B19").Value = "Verschillende oefenmodussen"

$backlog.Range("B3").Select()

# ---------------------------------------------------------------------------
# Sprint 2 sheet: flesh out the two placeholder rows (7-8) with the new user
# stories for score & lives tracking.
# ---------------------------------------------------------------------------
$sprint2 = $wb.Worksheets.Item("Sprint 2")

$sprint2.Range("A7").Value = 13
$sprint2.Range("B7").Value = "Als een user wil ik zien hoeveel ik goed heb of mijn score zien, zodat ik mijn voortgang kan bekijken en mij te motiveren om door te gaan."
$sprint2.Range("C7").Value = "M"
$sprint2.Range("D7").Value = 3
$sprint2.Range("E7").Value = "1. Schijf code voor score.`n2. Laat de score zien op het oefenscherm"
$sprint2.Rows.Item(7).RowHeight = 30

$sprint2.Range("A8").Value = 14
$sprint2.Range("B8").Value = "Als een user wil ik zien hoeveel levens / kansen ik nog heb."
$sprint2.Range("C8").Value = "M"
$sprint2.Range("D8").Value = 3
$sprint2.Range("E8").Value = "1. Schrijf code voor levens.`n2. Wanneer de gebruiker 0 levens over heeft, beëindig de sommen."
$sprint2.Range("F8").Value = " Wanneer de gebruiker 0 levens over heeft, beëindig de sommen. Laat de levens zien op het oefenscherm."
$sprint2.Rows.Item(8).RowHeight = 75

$sprint2.Range("E8").Select()
